# Update "想去人数" (F column) figures that were refreshed in the new
# data export (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 949
$ws1.Range("F5").Value  = 533
$ws1.Range("F6").Value  = 8035
$ws1.Range("F9").Value  = 1942
$ws1.Range("F10").Value = 5795
$ws1.Range("F12").Value = 239
$ws1.Range("F13").Value = 308
$ws1.Range("F14").Value = 8175
$ws1.Range("F15").Value = 9538
$ws1.Range("F16").Value = 1165
$ws1.Range("F17").Value = 954
$ws1.Range("F18").Value = 4591
$ws1.Range("F19").Value = 712
$ws1.Range("F20").Value = 288
$ws1.Range("F22").Value = 297
$ws1.Range("F25").Value = 138
$ws1.Range("F26").Value = 1738
$ws1.Range("F27").Value = 772
$ws1.Range("F28").Value = 1005
$ws1.Range("F29").Value = 212
$ws1.Range("F30").Value = 1925
$ws1.Range("F31").Value = 357
$ws1.Range("F32").Value = 498
$ws1.Range("F34").Value = 309
$ws1.Range("F36").Value = 1522
$ws1.Range("F38").Value = 1332
$ws1.Range("F39").Value = 14
$ws1.Range("F40").Value = 825
$ws1.Range("F48").Value = 187
$ws1.Range("F49").Value = 4127

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 0

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 949
$ws4.Range("F6").Value  = 533
$ws4.Range("F7").Value  = 8035
$ws4.Range("F11").Value = 5795
$ws4.Range("F13").Value = 8175
$ws4.Range("F14").Value = 9538
$ws4.Range("F16").Value = 1165
$ws4.Range("F17").Value = 954
$ws4.Range("F18").Value = 4591
$ws4.Range("F19").Value = 712
$ws4.Range("F20").Value = 288
$ws4.Range("F22").Value = 297
$ws4.Range("F26").Value = 138
$ws4.Range("F27").Value = 1738
$ws4.Range("F28").Value = 772
$ws4.Range("F29").Value = 1005
$ws4.Range("F30").Value = 212
$ws4.Range("F32").Value = 1925
$ws4.Range("F33").Value = 357
$ws4.Range("F38").Value = 825
$ws4.Range("F47").Value = 187
$ws4.Range("F48").Value = 4127
